$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("پدافند") values replaced with yes/no ("بله"/"خیر") instead of
# the old strength scale (زیاد/متوسط/کم), and column G ("نوع خاک") values
# remapped from قوی/متوسط/ضعیف to متراکم/سست. Also "کوهستان" -> "کوهستانی"
# in column H where it appears (rows 8 and 11), as per the commit
# "add dist as backup".

$ws.Range("C2").Value = "بله"
$ws.Range("G2").Value = "متراکم"

$ws.Range("C3").Value = "بله"
$ws.Range("G3").Value = "سست"

$ws.Range("C4").Value = "بله"
$ws.Range("G4").Value = "سست"

$ws.Range("C5").Value = "بله"
$ws.Range("G5").Value = "سست"

$ws.Range("C6").Value = "بله"
$ws.Range("G6").Value = "متراکم"

$ws.Range("C7").Value = "خیر"
$ws.Range("G7").Value = "سست"

$ws.Range("C8").Value = "خیر"
$ws.Range("G8").Value = "سست"
$ws.Range("H8").Value = "کوهستانی"

$ws.Range("C9").Value = "بله"
$ws.Range("G9").Value = "سست"

$ws.Range("C10").Value = "خیر"
$ws.Range("G10").Value = "متراکم"

$ws.Range("C11").Value = "بله"
$ws.Range("G11").Value = "سست"
$ws.Range("H11").Value = "کوهستانی"

# Update the active selection to match the saved state (I11).
$ws.Range("I11").Select()
